# Insert a new weekly price record for "Terminal La Palmera de La Serena - Coco".
# This shifts the existing rows 23-26 down to 24-27 (preserving their values)
# and inserts a brand-new row 23 with the latest week's data, extending the
# sheet's used range from A1:T26 to A1:T27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 23..26 down to 24..27, inheriting formatting from the row above
# (this is how Excel's own Insert behaves, matching style s="2" on column D).
$ws.Rows("23").Insert()

# Populate the newly inserted row 23 with the new record's values.
$ws.Range("A23").Value = 8
$ws.Range("B23").Value = "Terminal La Palmera de La Serena"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 44466
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100108
$ws.Range("H23").Value = "Tropicales y subtropicales"
$ws.Range("I23").Value = 100108007
$ws.Range("J23").Value = "Coco"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 20000
$ws.Range("O23").Value = 21000
$ws.Range("P23").Value = 20500
$ws.Range("Q23").Value = "$/malla 20 unidades"
$ws.Range("R23").Value = "Perú"
$ws.Range("S23").Value = 1025
$ws.Range("T23").Value = 20
